$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 9.486944736877136
$ws.Range("C2").Value = 5.692037375433699
$ws.Range("E2").Value = 16.46761802819981
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.625857913082736
$ws.Range("K2").Value = 8.689676238649927
$ws.Range("N2").Value = 17.83207224655341
$ws.Range("O2").Value = 20.88780721578135
$ws.Range("B3").Value = 9.176147000957604
$ws.Range("C3").Value = 5.533617818523892
$ws.Range("E3").Value = 15.53660390023953
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.627660703615569
$ws.Range("K3").Value = 8.47293520291586
$ws.Range("N3").Value = 17.89440262231351
$ws.Range("O3").Value = 20.96923847926976
$ws.Range("B4").Value = 8.981637274337329
$ws.Range("C4").Value = 5.433129413385758
$ws.Range("E4").Value = 14.94022716310238
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.628824678918332
$ws.Range("K4").Value = 8.338501143658688
$ws.Range("N4").Value = 17.93445351585281
$ws.Range("O4").Value = 21.02431268934053
$ws.Range("B5").Value = 8.901576866738212
$ws.Range("C5").Value = 5.391407216478807
$ws.Range("E5").Value = 14.69125042349741
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.629313403159999
$ws.Range("K5").Value = 8.283459738819911
$ws.Range("N5").Value = 17.95122357099906
$ws.Range("O5").Value = 21.04802752966808
$ws.Range("B6").Value = 8.888238780418796
$ws.Range("C6").Value = 5.384433765446996
$ws.Range("E6").Value = 14.64955770032874
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.629395426273974
$ws.Range("K6").Value = 8.274307017544775
$ws.Range("N6").Value = 17.95403538641524
$ws.Range("O6").Value = 21.05204203231929
$ws.Range("B7").Value = 8.980560597305573
$ws.Range("C7").Value = 5.432569809545079
$ws.Range("E7").Value = 14.93689306454285
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.628831211676858
$ws.Range("K7").Value = 8.337759771256398
$ws.Range("N7").Value = 17.93467786283122
$ws.Range("O7").Value = 21.02462737316353
$ws.Range("B8").Value = 9.380619351871719
$ws.Range("C8").Value = 5.638103435200865
$ws.Range("E8").Value = 16.15188364314047
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.626467702531456
$ws.Range("K8").Value = 8.615273058501405
$ws.Range("N8").Value = 17.85319523949643
$ws.Range("O8").Value = 20.91482868218441
$ws.Range("B9").Value = 10.13080589907931
$ws.Range("C9").Value = 6.014152415137067
$ws.Range("E9").Value = 18.40459232081288
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.622283374254512
$ws.Range("K9").Value = 9.145326943500081
$ws.Range("N9").Value = 17.70746293219693
$ws.Range("O9").Value = 20.73998225737555
$ws.Range("B10").Value = 10.65489516241233
$ws.Range("C10").Value = 6.272218268324738
$ws.Range("E10").Value = 20.03798457283578
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.619480704599982
$ws.Range("K10").Value = 9.521887977481956
$ws.Range("N10").Value = 17.60886659828725
$ws.Range("O10").Value = 20.63646353469302
$ws.Range("B11").Value = 10.88634380457577
$ws.Range("C11").Value = 6.38534824735527
$ws.Range("E11").Value = 20.73892322596665
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.618264005718688
$ws.Range("K11").Value = 9.689599124429529
$ws.Range("N11").Value = 17.56583216427281
$ws.Range("O11").Value = 20.59484241489709
$ws.Range("B12").Value = 10.9729080312953
$ws.Range("C12").Value = 6.42755138834333
$ws.Range("E12").Value = 20.99833245944162
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.61781159979962
$ws.Range("K12").Value = 9.752530915452578
$ws.Range("N12").Value = 17.54979599442214
$ws.Range("O12").Value = 20.57987271553377
$ws.Range("B13").Value = 10.9543140907918
$ws.Range("C13").Value = 6.418490851263345
$ws.Range("E13").Value = 20.94273130535957
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.617908663669867
$ws.Range("K13").Value = 9.739003995138537
$ws.Range("N13").Value = 17.55323812627335
$ws.Range("O13").Value = 20.58306143998414
$ws.Range("B14").Value = 10.8934875691264
$ws.Range("C14").Value = 6.388833211173075
$ws.Range("E14").Value = 20.76038546259011
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.618226619289415
$ws.Range("K14").Value = 9.694788450587005
$ws.Range("N14").Value = 17.5645076549495
$ws.Range("O14").Value = 20.593594962874
$ws.Range("B15").Value = 10.85608667808933
$ws.Range("C15").Value = 6.370583478466188
$ws.Range("E15").Value = 20.64791022520905
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.618422459978841
$ws.Range("K15").Value = 9.667628320889365
$ws.Range("N15").Value = 17.57144439200766
$ws.Range("O15").Value = 20.60015024095437
$ws.Range("B16").Value = 10.63962237480754
$ws.Range("C16").Value = 6.26473707265304
$ws.Range("E16").Value = 19.99133224459941
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.619561387038201
$ws.Range("K16").Value = 9.510850142265079
$ws.Range("N16").Value = 17.61171545983826
$ws.Range("O16").Value = 20.63929405079681
$ws.Range("B17").Value = 10.50498609809082
$ws.Range("C17").Value = 6.198693829943362
$ws.Range("E17").Value = 19.57778363809619
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.620274969140302
$ws.Range("K17").Value = 9.413707900273948
$ws.Range("N17").Value = 17.63688502464709
$ws.Range("O17").Value = 20.66471191998896
$ws.Range("B18").Value = 10.42689569378371
$ws.Range("C18").Value = 6.160307250461537
$ws.Range("E18").Value = 19.33595684735339
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.620690888071648
$ws.Range("K18").Value = 9.357499811571058
$ws.Range("N18").Value = 17.65153304567448
$ws.Range("O18").Value = 20.67984619781825
$ws.Range("B19").Value = 10.40034634661504
$ws.Range("C19").Value = 6.147242222823314
$ws.Range("E19").Value = 19.25339652727322
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.620832654560925
$ws.Range("K19").Value = 9.33841321250744
$ws.Range("N19").Value = 17.65652205012461
$ws.Range("O19").Value = 20.68505864029381
$ws.Range("B20").Value = 10.51938635656952
$ws.Range("C20").Value = 6.205765847620163
$ws.Range("E20").Value = 19.62221667262738
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.620198439715091
$ws.Range("K20").Value = 9.424083927305656
$ws.Range("N20").Value = 17.63418797739995
$ws.Range("O20").Value = 20.66195285149945
$ws.Range("B21").Value = 10.91138369999933
$ws.Range("C21").Value = 6.397561833529611
$ws.Range("E21").Value = 20.81410800481927
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.618133002256005
$ws.Range("K21").Value = 9.707791740970906
$ws.Range("N21").Value = 17.56119047437193
$ws.Range("O21").Value = 20.59047949862632
$ws.Range("B22").Value = 11.16124899018282
$ws.Range("C22").Value = 6.519191133076338
$ws.Range("E22").Value = 21.55799785061672
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.616831661544797
$ws.Range("K22").Value = 9.889824608193329
$ws.Range("N22").Value = 17.51499747134682
$ws.Range("O22").Value = 20.54838204156139
$ws.Range("B23").Value = 11.02849403083878
$ws.Range("C23").Value = 6.454622798326046
$ws.Range("E23").Value = 21.16416886920414
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.617521784662087
$ws.Range("K23").Value = 9.792998901345367
$ws.Range("N23").Value = 17.53951335014584
$ws.Range("O23").Value = 20.57042648392294
$ws.Range("B24").Value = 10.51287813098489
$ws.Range("C24").Value = 6.202569886924591
$ws.Range("E24").Value = 19.60214118473799
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.620233020998713
$ws.Range("K24").Value = 9.419394041229115
$ws.Range("N24").Value = 17.63540675878595
$ws.Range("O24").Value = 20.66319860311318
$ws.Range("B25").Value = 9.932210334012236
$ws.Range("C25").Value = 5.915504341227162
$ws.Range("E25").Value = 17.765529280945
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.623367436024327
$ws.Range("K25").Value = 9.003904660987823
$ws.Range("N25").Value = 17.74539257816797
$ws.Range("O25").Value = 20.78292222778792
